# Update the arithmetic expressions in the practice worksheet table.
# Each cell's text is unique within the document, so MatchWholeWord
# Find/Replace (Replace All) safely targets exactly one cell each time.
$d = $word.ActiveDocument

$d.Content.Find.Execute("76-49=", $true, $false, $false, $false, $false, $true, 1, $false, "61+13=", 2) | Out-Null
$d.Content.Find.Execute("12+45=", $true, $false, $false, $false, $false, $true, 1, $false, "20+56=", 2) | Out-Null
$d.Content.Find.Execute("55-35=", $true, $false, $false, $false, $false, $true, 1, $false, "92-52=", 2) | Out-Null
$d.Content.Find.Execute("25+73=", $true, $false, $false, $false, $false, $true, 1, $false, "68-32=", 2) | Out-Null
$d.Content.Find.Execute("12+64=", $true, $false, $false, $false, $false, $true, 1, $false, "91-42=", 2) | Out-Null
$d.Content.Find.Execute("87-42=", $true, $false, $false, $false, $false, $true, 1, $false, "5+73=", 2) | Out-Null
$d.Content.Find.Execute("48-27=", $true, $false, $false, $false, $false, $true, 1, $false, "38+37=", 2) | Out-Null
$d.Content.Find.Execute("65+2=", $true, $false, $false, $false, $false, $true, 1, $false, "18-1=", 2) | Out-Null
$d.Content.Find.Execute("76-10=", $true, $false, $false, $false, $false, $true, 1, $false, "96-21=", 2) | Out-Null
$d.Content.Find.Execute("59-12=", $true, $false, $false, $false, $false, $true, 1, $false, "64-19=", 2) | Out-Null
$d.Content.Find.Execute("65+11=", $true, $false, $false, $false, $false, $true, 1, $false, "34+36=", 2) | Out-Null
$d.Content.Find.Execute("88-69=", $true, $false, $false, $false, $false, $true, 1, $false, "70-18=", 2) | Out-Null
$d.Content.Find.Execute("61-25=", $true, $false, $false, $false, $false, $true, 1, $false, "94-81=", 2) | Out-Null
$d.Content.Find.Execute("97-63=", $true, $false, $false, $false, $false, $true, 1, $false, "71-36=", 2) | Out-Null
$d.Content.Find.Execute("69-67=", $true, $false, $false, $false, $false, $true, 1, $false, "92-42=", 2) | Out-Null
$d.Content.Find.Execute("73-45=", $true, $false, $false, $false, $false, $true, 1, $false, "46+41=", 2) | Out-Null
$d.Content.Find.Execute("59-26=", $true, $false, $false, $false, $false, $true, 1, $false, "2+16=", 2) | Out-Null
$d.Content.Find.Execute("16-4=", $true, $false, $false, $false, $false, $true, 1, $false, "64+6=", 2) | Out-Null
$d.Content.Find.Execute("60+26=", $true, $false, $false, $false, $false, $true, 1, $false, "21+38=", 2) | Out-Null
$d.Content.Find.Execute("81+14=", $true, $false, $false, $false, $false, $true, 1, $false, "82-26=", 2) | Out-Null
$d.Content.Find.Execute("45-32=", $true, $false, $false, $false, $false, $true, 1, $false, "71-36=", 2) | Out-Null
$d.Content.Find.Execute("33+66=", $true, $false, $false, $false, $false, $true, 1, $false, "46+49=", 2) | Out-Null
$d.Content.Find.Execute("40+36=", $true, $false, $false, $false, $false, $true, 1, $false, "66+23=", 2) | Out-Null
$d.Content.Find.Execute("62-61=", $true, $false, $false, $false, $false, $true, 1, $false, "85-72=", 2) | Out-Null
$d.Content.Find.Execute("46-25=", $true, $false, $false, $false, $false, $true, 1, $false, "88-58=", 2) | Out-Null
$d.Content.Find.Execute("74-20=", $true, $false, $false, $false, $false, $true, 1, $false, "83+4=", 2) | Out-Null
$d.Content.Find.Execute("49-31=", $true, $false, $false, $false, $false, $true, 1, $false, "1+12=", 2) | Out-Null
$d.Content.Find.Execute("84-60=", $true, $false, $false, $false, $false, $true, 1, $false, "8+12=", 2) | Out-Null
$d.Content.Find.Execute("22+42=", $true, $false, $false, $false, $false, $true, 1, $false, "64-51=", 2) | Out-Null
$d.Content.Find.Execute("78-16=", $true, $false, $false, $false, $false, $true, 1, $false, "21-4=", 2) | Out-Null
$d.Content.Find.Execute("54-14=", $true, $false, $false, $false, $false, $true, 1, $false, "5+92=", 2) | Out-Null
$d.Content.Find.Execute("3+6=", $true, $false, $false, $false, $false, $true, 1, $false, "68-12=", 2) | Out-Null
$d.Content.Find.Execute("61-47=", $true, $false, $false, $false, $false, $true, 1, $false, "62-62=", 2) | Out-Null
$d.Content.Find.Execute("63-47=", $true, $false, $false, $false, $false, $true, 1, $false, "79-37=", 2) | Out-Null
$d.Content.Find.Execute("73-12=", $true, $false, $false, $false, $false, $true, 1, $false, "35+0=", 2) | Out-Null
$d.Content.Find.Execute("90-82=", $true, $false, $false, $false, $false, $true, 1, $false, "52-39=", 2) | Out-Null
$d.Content.Find.Execute("29-2=", $true, $false, $false, $false, $false, $true, 1, $false, "27+68=", 2) | Out-Null
$d.Content.Find.Execute("68-59=", $true, $false, $false, $false, $false, $true, 1, $false, "76+18=", 2) | Out-Null
$d.Content.Find.Execute("86-73=", $true, $false, $false, $false, $false, $true, 1, $false, "39+31=", 2) | Out-Null
$d.Content.Find.Execute("59-22=", $true, $false, $false, $false, $false, $true, 1, $false, "8+52=", 2) | Out-Null
$d.Content.Find.Execute("73+5=", $true, $false, $false, $false, $false, $true, 1, $false, "68-33=", 2) | Out-Null
$d.Content.Find.Execute("22+5=", $true, $false, $false, $false, $false, $true, 1, $false, "82-53=", 2) | Out-Null
$d.Content.Find.Execute("71+13=", $true, $false, $false, $false, $false, $true, 1, $false, "82+16=", 2) | Out-Null
$d.Content.Find.Execute("38+4=", $true, $false, $false, $false, $false, $true, 1, $false, "64+2=", 2) | Out-Null
$d.Content.Find.Execute("63-37=", $true, $false, $false, $false, $false, $true, 1, $false, "67-0=", 2) | Out-Null
$d.Content.Find.Execute("52+26=", $true, $false, $false, $false, $false, $true, 1, $false, "59+32=", 2) | Out-Null
$d.Content.Find.Execute("11+79=", $true, $false, $false, $false, $false, $true, 1, $false, "77-34=", 2) | Out-Null
$d.Content.Find.Execute("63-55=", $true, $false, $false, $false, $false, $true, 1, $false, "7+35=", 2) | Out-Null
$d.Content.Find.Execute("51+21=", $true, $false, $false, $false, $false, $true, 1, $false, "45+11=", 2) | Out-Null
$d.Content.Find.Execute("94-55=", $true, $false, $false, $false, $false, $true, 1, $false, "55+7=", 2) | Out-Null
$d.Content.Find.Execute("49-20=", $true, $false, $false, $false, $false, $true, 1, $false, "98-21=", 2) | Out-Null
$d.Content.Find.Execute("62+30=", $true, $false, $false, $false, $false, $true, 1, $false, "7+21=", 2) | Out-Null
$d.Content.Find.Execute("87-87=", $true, $false, $false, $false, $false, $true, 1, $false, "99-4=", 2) | Out-Null
$d.Content.Find.Execute("52+37=", $true, $false, $false, $false, $false, $true, 1, $false, "65-10=", 2) | Out-Null
$d.Content.Find.Execute("57+33=", $true, $false, $false, $false, $false, $true, 1, $false, "52+40=", 2) | Out-Null
$d.Content.Find.Execute("2+85=", $true, $false, $false, $false, $false, $true, 1, $false, "18+77=", 2) | Out-Null
$d.Content.Find.Execute("62+15=", $true, $false, $false, $false, $false, $true, 1, $false, "16+33=", 2) | Out-Null
$d.Content.Find.Execute("50+40=", $true, $false, $false, $false, $false, $true, 1, $false, "63-28=", 2) | Out-Null
$d.Content.Find.Execute("15+54=", $true, $false, $false, $false, $false, $true, 1, $false, "92-40=", 2) | Out-Null
$d.Content.Find.Execute("2+95=", $true, $false, $false, $false, $false, $true, 1, $false, "86-72=", 2) | Out-Null
$d.Content.Find.Execute("97-30=", $true, $false, $false, $false, $false, $true, 1, $false, "22+37=", 2) | Out-Null
$d.Content.Find.Execute("32+1=", $true, $false, $false, $false, $false, $true, 1, $false, "93-29=", 2) | Out-Null
$d.Content.Find.Execute("20+22=", $true, $false, $false, $false, $false, $true, 1, $false, "19+12=", 2) | Out-Null
$d.Content.Find.Execute("39+2=", $true, $false, $false, $false, $false, $true, 1, $false, "67+2=", 2) | Out-Null
$d.Content.Find.Execute("43-20=", $true, $false, $false, $false, $false, $true, 1, $false, "67-26=", 2) | Out-Null
$d.Content.Find.Execute("35+34=", $true, $false, $false, $false, $false, $true, 1, $false, "70-51=", 2) | Out-Null
$d.Content.Find.Execute("81-76=", $true, $false, $false, $false, $false, $true, 1, $false, "9-8=", 2) | Out-Null
$d.Content.Find.Execute("96-81=", $true, $false, $false, $false, $false, $true, 1, $false, "96-86=", 2) | Out-Null
$d.Content.Find.Execute("93-72=", $true, $false, $false, $false, $false, $true, 1, $false, "52-50=", 2) | Out-Null
$d.Content.Find.Execute("29+2=", $true, $false, $false, $false, $false, $true, 1, $false, "61-36=", 2) | Out-Null
$d.Content.Find.Execute("34+57=", $true, $false, $false, $false, $false, $true, 1, $false, "87-71=", 2) | Out-Null
$d.Content.Find.Execute("46-1=", $true, $false, $false, $false, $false, $true, 1, $false, "78+18=", 2) | Out-Null
$d.Content.Find.Execute("40-27=", $true, $false, $false, $false, $false, $true, 1, $false, "70+9=", 2) | Out-Null
$d.Content.Find.Execute("54+28=", $true, $false, $false, $false, $false, $true, 1, $false, "9+46=", 2) | Out-Null
$d.Content.Find.Execute("3+55=", $true, $false, $false, $false, $false, $true, 1, $false, "95-79=", 2) | Out-Null
$d.Content.Find.Execute("81+5=", $true, $false, $false, $false, $false, $true, 1, $false, "79-7=", 2) | Out-Null
$d.Content.Find.Execute("84-61=", $true, $false, $false, $false, $false, $true, 1, $false, "28-27=", 2) | Out-Null
$d.Content.Find.Execute("44+53=", $true, $false, $false, $false, $false, $true, 1, $false, "75-73=", 2) | Out-Null
$d.Content.Find.Execute("47-33=", $true, $false, $false, $false, $false, $true, 1, $false, "80-39=", 2) | Out-Null
$d.Content.Find.Execute("71+14=", $true, $false, $false, $false, $false, $true, 1, $false, "75-58=", 2) | Out-Null
$d.Content.Find.Execute("58+20=", $true, $false, $false, $false, $false, $true, 1, $false, "99-41=", 2) | Out-Null
$d.Content.Find.Execute("14+47=", $true, $false, $false, $false, $false, $true, 1, $false, "46-16=", 2) | Out-Null
$d.Content.Find.Execute("1+85=", $true, $false, $false, $false, $false, $true, 1, $false, "70+26=", 2) | Out-Null
$d.Content.Find.Execute("66-62=", $true, $false, $false, $false, $false, $true, 1, $false, "17+51=", 2) | Out-Null
$d.Content.Find.Execute("17+46=", $true, $false, $false, $false, $false, $true, 1, $false, "53+29=", 2) | Out-Null
$d.Content.Find.Execute("86-58=", $true, $false, $false, $false, $false, $true, 1, $false, "72-38=", 2) | Out-Null
$d.Content.Find.Execute("91-46=", $true, $false, $false, $false, $false, $true, 1, $false, "35-32=", 2) | Out-Null
$d.Content.Find.Execute("77-55=", $true, $false, $false, $false, $false, $true, 1, $false, "26+64=", 2) | Out-Null
$d.Content.Find.Execute("99-48=", $true, $false, $false, $false, $false, $true, 1, $false, "74-17=", 2) | Out-Null
$d.Content.Find.Execute("19+59=", $true, $false, $false, $false, $false, $true, 1, $false, "72-24=", 2) | Out-Null
$d.Content.Find.Execute("19+62=", $true, $false, $false, $false, $false, $true, 1, $false, "2+22=", 2) | Out-Null
$d.Content.Find.Execute("16+6=", $true, $false, $false, $false, $false, $true, 1, $false, "7+66=", 2) | Out-Null
$d.Content.Find.Execute("82-46=", $true, $false, $false, $false, $false, $true, 1, $false, "78-47=", 2) | Out-Null
$d.Content.Find.Execute("65-25=", $true, $false, $false, $false, $false, $true, 1, $false, "82-57=", 2) | Out-Null
$d.Content.Find.Execute("91-64=", $true, $false, $false, $false, $false, $true, 1, $false, "56-19=", 2) | Out-Null
$d.Content.Find.Execute("23+72=", $true, $false, $false, $false, $false, $true, 1, $false, "81-11=", 2) | Out-Null
$d.Content.Find.Execute("32+32=", $true, $false, $false, $false, $false, $true, 1, $false, "85+10=", 2) | Out-Null
$d.Content.Find.Execute("90-73=", $true, $false, $false, $false, $false, $true, 1, $false, "14+66=", 2) | Out-Null
$d.Content.Find.Execute("85-13=", $true, $false, $false, $false, $false, $true, 1, $false, "91-69=", 2) | Out-Null
$d.Content.Find.Execute("70-30=", $true, $false, $false, $false, $false, $true, 1, $false, "96-34=", 2) | Out-Null
